# Applies the "Updated symbol list" edit described in the diff:
#  - Rows 6-17: coin Name/Link shift up one slot (a new coin,
#    GateToken, is inserted at row 6 and the remaining rows shift
#    down), plus fresh Price/Volume(1h) figures throughout.
#  - Several other rows only get updated Price/Volume(1h) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a literal
# TEXT value (the sheet stores Price/Volume(1h) as text, e.g. "302.77"
# or "2.12%" -- a plain Value2 assignment would make Excel parse these
# numeric-/percent-looking strings into real numbers). We build the
# text via a throwaway formula (="...") and then collapse the cell to
# a value-only paste of itself, which keeps the result a plain string
# without touching the cells style/number format.
function Set-TextValue {
    param($Ref, $Text)
    $cell = $ws.Range($Ref)
    $escaped = $Text.Replace("""", """""")
    $cell.Formula = "=""" + $escaped + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Coin (B) and Link (C) columns: plain text, safe to assign directly ---
$ws.Range('B6').Value2 = 'GateToken'
$ws.Range('C6').Value2 = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B7').Value2 = 'FTXToken'
$ws.Range('C7').Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('B8').Value2 = 'MXToken'
$ws.Range('C8').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B9').Value2 = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value2 = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B10').Value2 = 'WazirX'
$ws.Range('C10').Value2 = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B11').Value2 = 'MandalaExchangeToken'
$ws.Range('C11').Value2 = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B12').Value2 = 'BitrueCoin'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B13').Value2 = 'BitMartToken'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B14').Value2 = 'BitForexToken'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B15').Value2 = 'TigerCash'
$ws.Range('C15').Value2 = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B16').Value2 = 'UpBots'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('B17').Value2 = 'LEO'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

# --- Price (D) and Volume(1h) (E) columns: numeric-looking text, must
#     stay literal strings ---
Set-TextValue 'D2' '302.77'
Set-TextValue 'E2' '2.12%'
Set-TextValue 'D3' '43.29'
Set-TextValue 'E3' '4.86%'
Set-TextValue 'D4' '5.066'
Set-TextValue 'E4' '0.40%'
Set-TextValue 'D5' '0.07689'
Set-TextValue 'E5' '3.22%'
Set-TextValue 'D6' '4.420'
Set-TextValue 'E6' '1.50%'
Set-TextValue 'D7' '1.621'
Set-TextValue 'E7' '3.00%'
Set-TextValue 'D8' '1.045'
Set-TextValue 'E8' '12.11%'
Set-TextValue 'D9' '0.1260'
Set-TextValue 'E9' '5.18%'
Set-TextValue 'D10' '0.1870'
Set-TextValue 'E10' '3.86%'
Set-TextValue 'D11' '0.09103'
Set-TextValue 'E11' '3.42%'
Set-TextValue 'D12' '0.04167'
Set-TextValue 'E12' '-3.46%'
Set-TextValue 'D13' '0.1049'
Set-TextValue 'E13' '0.29%'
Set-TextValue 'D14' '0.001294'
Set-TextValue 'E14' '1.38%'
Set-TextValue 'D15' '0.005753'
Set-TextValue 'E15' '-1.76%'
Set-TextValue 'D16' '0.007430'
Set-TextValue 'E16' '1,897.29%'
Set-TextValue 'D17' '3.345'
Set-TextValue 'E17' '-0.40%'
Set-TextValue 'E18' '-1.95%'
Set-TextValue 'D19' '0.3353'
Set-TextValue 'E19' '1.36%'
Set-TextValue 'D20' '8.662'
Set-TextValue 'E20' '8.10%'
Set-TextValue 'E21' '-0.85%'
Set-TextValue 'D23' '0.04165'
Set-TextValue 'E23' '3.92%'
Set-TextValue 'D24' '0.001284'
Set-TextValue 'E24' '1.50%'
Set-TextValue 'D25' '0.004453'
Set-TextValue 'E25' '15.32%'
Set-TextValue 'D26' '0.0001349'
Set-TextValue 'E26' '9.82%'
Set-TextValue 'D38' '0.02457'
Set-TextValue 'E38' '3.58%'
Set-TextValue 'D39' '0.05278'
Set-TextValue 'E39' '2.36%'
Set-TextValue 'D40' '0.005931'
Set-TextValue 'E40' '-1.66%'
Set-TextValue 'D41' '0.007696'
Set-TextValue 'E41' '-1.22%'
Set-TextValue 'D42' '0.1346'
Set-TextValue 'E42' '2.50%'
Set-TextValue 'D43' '0.007375'
Set-TextValue 'E43' '-0.30%'
Set-TextValue 'D44' '0.007566'
Set-TextValue 'E44' '-3.20%'
Set-TextValue 'D45' '0.3009'
Set-TextValue 'E45' '2.33%'
Set-TextValue 'D46' '0.00006699'
Set-TextValue 'E46' '4.18%'
Set-TextValue 'E47' '-0.05%'
Set-TextValue 'E48' '-29.27%'
Set-TextValue 'E49' '0.06%'
Set-TextValue 'D50' '0.00002098'
Set-TextValue 'E50' '-0.05%'
Set-TextValue 'D51' '0.0001998'
Set-TextValue 'E51' '-0.05%'

$wb.Application.CutCopyMode = $false

